$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, "H").Value = 846.0345
$ws.Cells.Item(19, "I").Value = 290.83334
$ws.Cells.Item(19, "K").Value = 290.83334
$ws.Cells.Item(19, "M").Value = -115.83334
$ws.Cells.Item(32, "H").Value = 1969.8
$ws.Cells.Item(32, "I").Value = 1949.6666
$ws.Cells.Item(32, "J").Value = 2000
$ws.Cells.Item(32, "K").Value = 1949.6666
$ws.Cells.Item(32, "L").Value = 2000
$ws.Cells.Item(32, "M").Value = -1623.6666
$ws.Cells.Item(32, "N").Value = -2652
$ws.Cells.Item(55, "H").Value = 362.58334
$ws.Cells.Item(55, "I").Value = 353
$ws.Cells.Item(55, "J").Value = 376
$ws.Cells.Item(55, "K").Value = 353
$ws.Cells.Item(55, "L").Value = 376
$ws.Cells.Item(55, "M").Value = -139
$ws.Cells.Item(55, "N").Value = -804
$ws.Cells.Item(82, "H").Value = 546
$ws.Cells.Item(82, "J").Value = 0
$ws.Cells.Item(82, "L").Value = 0
$ws.Cells.Item(82, "N").ClearContents()
$ws.Cells.Item(85, "H").Value = 546
$ws.Cells.Item(85, "J").Value = 0
$ws.Cells.Item(85, "L").Value = 0
$ws.Cells.Item(85, "N").ClearContents()
$ws.Cells.Item(116, "H").Value = 12497.5
$ws.Cells.Item(116, "I").Value = 14996.667
$ws.Cells.Item(116, "J").Value = 5000
$ws.Cells.Item(116, "K").Value = 14996.667
$ws.Cells.Item(116, "L").Value = 5000
$ws.Cells.Item(116, "M").Value = -11554.667
$ws.Cells.Item(116, "N").Value = -11884
$ws.Cells.Item(138, "H").Value = 5683541
$ws.Cells.Item(138, "I").Value = 1280.1
$ws.Cells.Item(138, "J").Value = 10418758
$ws.Cells.Item(138, "K").Value = 3840.3
$ws.Cells.Item(138, "L").Value = 31256274
$ws.Cells.Item(138, "M").Value = 1299.7
$ws.Cells.Item(138, "N").Value = -31266554

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, "H").Value = 2231.8
$ws.Cells.Item(45, "I").Value = 2146.4443
$ws.Cells.Item(45, "K").Value = 2146.4443
$ws.Cells.Item(45, "M").Value = -1769.4443
$ws.Cells.Item(122, "H").Value = 8691.352999999999
$ws.Cells.Item(122, "I").Value = 9845.643
$ws.Cells.Item(122, "J").Value = 3304.6667
$ws.Cells.Item(122, "K").Value = 29536.929
$ws.Cells.Item(122, "L").Value = 9914.000100000001
$ws.Cells.Item(122, "M").Value = -27086.929
$ws.Cells.Item(122, "N").Value = -14814.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, "H").Value = 4541.2905
$ws.Cells.Item(105, "I").Value = 2930
$ws.Cells.Item(105, "J").Value = 4928
$ws.Cells.Item(105, "K").Value = 2930
$ws.Cells.Item(105, "L").Value = 4928
$ws.Cells.Item(105, "M").Value = -1183
$ws.Cells.Item(105, "N").Value = -8422
$ws.Cells.Item(134, "H").Value = 3559.6667
$ws.Cells.Item(134, "I").Value = 2419
$ws.Cells.Item(134, "J").Value = 6981.6665
$ws.Cells.Item(134, "K").Value = 7257
$ws.Cells.Item(134, "L").Value = 20944.9995
$ws.Cells.Item(134, "M").Value = -4722
$ws.Cells.Item(134, "N").Value = -26014.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, "H").Value = 3535.28
$ws.Cells.Item(94, "I").Value = 2339.2
$ws.Cells.Item(94, "J").Value = 4332.6665
$ws.Cells.Item(94, "K").Value = 2339.2
$ws.Cells.Item(94, "L").Value = 4332.6665
$ws.Cells.Item(94, "M").Value = -1888.2
$ws.Cells.Item(94, "N").Value = -5234.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, "H").Value = 94.2
$ws.Cells.Item(23, "I").Value = 58.4
$ws.Cells.Item(23, "K").Value = 175.2
$ws.Cells.Item(23, "M").Value = 59.80000000000001
$ws.Cells.Item(68, "H").Value = 1134.3334
$ws.Cells.Item(68, "I").Value = 501
$ws.Cells.Item(68, "K").Value = 1503
$ws.Cells.Item(68, "M").Value = -692
$ws.Cells.Item(71, "H").Value = 1134.3334
$ws.Cells.Item(71, "I").Value = 501
$ws.Cells.Item(71, "K").Value = 4509
$ws.Cells.Item(71, "M").Value = -453
$ws.Cells.Item(98, "H").Value = 310
$ws.Cells.Item(98, "J").Value = 198
$ws.Cells.Item(98, "L").Value = 594
$ws.Cells.Item(98, "N").Value = -3590
$ws.Cells.Item(118, "H").Value = 1966.0741
$ws.Cells.Item(118, "J").Value = 1831.6086
$ws.Cells.Item(118, "L").Value = 5494.825800000001
$ws.Cells.Item(118, "N").Value = -7980.825800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, "H").Value = 1571.7
$ws.Cells.Item(97, "I").Value = 1376.1538
$ws.Cells.Item(97, "K").Value = 1376.1538
$ws.Cells.Item(97, "M").Value = -880.1538
$ws.Cells.Item(113, "H").Value = 201102.2
$ws.Cells.Item(113, "I").Value = 201102.2
$ws.Cells.Item(113, "J").Value = 0
$ws.Cells.Item(113, "K").Value = 201102.2
$ws.Cells.Item(113, "L").Value = 0
$ws.Cells.Item(113, "M").Value = -198932.2
$ws.Cells.Item(113, "N").ClearContents()
$ws.Cells.Item(132, "H").Value = 4468.3076
$ws.Cells.Item(132, "I").Value = 3091.6667
$ws.Cells.Item(132, "J").Value = 7565.75
$ws.Cells.Item(132, "K").Value = 9275.000100000001
$ws.Cells.Item(132, "L").Value = 22697.25
$ws.Cells.Item(132, "M").Value = -6745.000100000001
$ws.Cells.Item(132, "N").Value = -27757.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, "H").Value = 6121.2104
$ws.Cells.Item(7, "I").Value = 6540.3
$ws.Cells.Item(7, "J").Value = 5655.5557
$ws.Cells.Item(7, "K").Value = 6540.3
$ws.Cells.Item(7, "L").Value = 5655.5557
$ws.Cells.Item(7, "M").Value = -6428.3
$ws.Cells.Item(7, "N").Value = -5879.5557
$ws.Cells.Item(22, "H").Value = 1550.2354
$ws.Cells.Item(22, "I").Value = 1220
$ws.Cells.Item(22, "J").Value = 1651.8462
$ws.Cells.Item(22, "K").Value = 1220
$ws.Cells.Item(22, "L").Value = 1651.8462
$ws.Cells.Item(22, "M").Value = -925
$ws.Cells.Item(22, "N").Value = -2241.8462
$ws.Cells.Item(27, "H").Value = 1550.2354
$ws.Cells.Item(27, "I").Value = 1220
$ws.Cells.Item(27, "J").Value = 1651.8462
$ws.Cells.Item(27, "K").Value = 1220
$ws.Cells.Item(27, "L").Value = 1651.8462
$ws.Cells.Item(27, "M").Value = -1113
$ws.Cells.Item(27, "N").Value = -1865.8462
$ws.Cells.Item(40, "H").Value = 4402.476
$ws.Cells.Item(40, "I").Value = 5540.4546
$ws.Cells.Item(40, "J").Value = 3150.7
$ws.Cells.Item(40, "K").Value = 5540.4546
$ws.Cells.Item(40, "L").Value = 3150.7
$ws.Cells.Item(40, "M").Value = -5404.4546
$ws.Cells.Item(40, "N").Value = -3422.7
$ws.Cells.Item(122, "H").Value = 8009.1113
$ws.Cells.Item(122, "I").Value = 9933.125
$ws.Cells.Item(122, "K").Value = 29799.375
$ws.Cells.Item(122, "M").Value = -27349.375
$ws.Cells.Item(126, "H").Value = 6121.2104
$ws.Cells.Item(126, "I").Value = 6540.3
$ws.Cells.Item(126, "J").Value = 5655.5557
$ws.Cells.Item(126, "K").Value = 19620.9
$ws.Cells.Item(126, "L").Value = 16966.6671
$ws.Cells.Item(126, "M").Value = -17150.9
$ws.Cells.Item(126, "N").Value = -21906.6671
$ws.Cells.Item(132, "H").Value = 10211943
$ws.Cells.Item(132, "I").Value = 5016.1333
$ws.Cells.Item(132, "J").Value = 26328142
$ws.Cells.Item(132, "K").Value = 15048.3999
$ws.Cells.Item(132, "L").Value = 78984426
$ws.Cells.Item(132, "M").Value = -12518.3999
$ws.Cells.Item(132, "N").Value = -78989486
$ws.Cells.Item(136, "H").Value = 17248940
$ws.Cells.Item(136, "I").Value = 21742474
$ws.Cells.Item(136, "J").Value = 23730.834
$ws.Cells.Item(136, "K").Value = 65227422
$ws.Cells.Item(136, "L").Value = 71192.50199999999
$ws.Cells.Item(136, "M").Value = -65224872
$ws.Cells.Item(136, "N").Value = -76292.50199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(120, "H").Value = 28000
$ws.Cells.Item(120, "I").Value = 25000
$ws.Cells.Item(120, "J").Value = 29500
$ws.Cells.Item(120, "K").Value = 25000
$ws.Cells.Item(120, "L").Value = 29500
$ws.Cells.Item(120, "M").Value = -20162
$ws.Cells.Item(120, "N").Value = -39176
$ws.Cells.Item(122, "H").Value = 1982.2903
$ws.Cells.Item(122, "I").Value = 1922.64
$ws.Cells.Item(122, "J").Value = 2230.8333
$ws.Cells.Item(122, "K").Value = 5767.92
$ws.Cells.Item(122, "L").Value = 6692.499899999999
$ws.Cells.Item(122, "M").Value = -3317.92
$ws.Cells.Item(122, "N").Value = -11592.4999
$ws.Cells.Item(126, "H").Value = 2936.5833
$ws.Cells.Item(126, "I").Value = 2132.389
$ws.Cells.Item(126, "J").Value = 5349.1665
$ws.Cells.Item(126, "K").Value = 6397.167
$ws.Cells.Item(126, "L").Value = 16047.4995
$ws.Cells.Item(126, "M").Value = -3927.167
$ws.Cells.Item(126, "N").Value = -20987.4995
$ws.Cells.Item(132, "H").Value = 1243.4324
$ws.Cells.Item(132, "I").Value = 913.13635
$ws.Cells.Item(132, "J").Value = 1727.8667
$ws.Cells.Item(132, "K").Value = 2739.40905
$ws.Cells.Item(132, "L").Value = 5183.6001
$ws.Cells.Item(132, "M").Value = -209.4090500000002
$ws.Cells.Item(132, "N").Value = -10243.6001

Write-Host "Applied all cell updates: $(188) sets, $(3) clears"
